$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - website names for each test
$ws.Range("C2").Value = "Web: GRAN WALKING"
$ws.Range("F2").Value = "Web: COOK IT"

# Row 5 - SEXO
$ws.Range("C5").Value = "Hombre"
$ws.Range("E5").Value = "Mujer"
$ws.Range("F5").Value = "Hombre"
$ws.Range("G5").Value = "Mujer"

# Row 6 - EDAD
$ws.Range("C6").Value = 23
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 31
$ws.Range("G6").Value = 29

# Row 7 - OCUPACION
$ws.Range("C7").Value = "Estudiante"
$ws.Range("E7").Value = "Profesora escolar"
$ws.Range("F7").Value = "Empresario"
$ws.Range("G7").Value = "Peluquera"

# Row 8 - EXPERIENCIA TIC
$ws.Range("C8").Value = "Experiencia en campos de programacion y html"
$ws.Range("E8").Value = "no"
$ws.Range("F8").Value = "Experiencia en diseño y analisis de webs"
$ws.Range("G8").Value = "Experiencia en creacion de paginas web usando aplicaciones automatizadas"

# Row 9 - PERFIL (describir)
$ws.Range("C9").Value = "Me gusta la programacion, estoy estudiando un modulo de informatica, tambien me gusta conocer gente nueva y explorar nuevas ciudades"
$ws.Range("E9").Value = "Me encantan los niños y la gente, soy una persona muy social y comprensiva"
$ws.Range("F9").Value = "Gestiono mis propios negocios y activos, me gusta valorar los buenos trabajos y la gastronomía en difrentes lugares"
$ws.Range("G9").Value = "Especialista en gastronomia, prueba y cata de alimentos"

# Rows 14-23 - SUS questionnaire answers (1-5) for the four users
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 5

$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1

$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 4

$ws.Range("C17").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 3

$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 5

$ws.Range("C19").Value = 3
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 1

$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3

$ws.Range("C21").Value = 2
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 3

$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4

$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2
